$wb = $excel.ActiveWorkbook

# @@ -1171,22 +1171,22 @@  (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 3665.6667
$ws.Range("I11").Value = 3665.6667
$ws.Range("K11").Value = 3665.6667
$ws.Range("M11").Value = -3525.6667

# @@ -1566,25 +1566,25 @@  (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1613.3636
$ws.Range("J19").Value = 1014.5
$ws.Range("L19").Value = 1014.5
$ws.Range("N19").Value = -1364.5

# @@ -5628,25 +5628,25 @@  (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 14789.111
$ws.Range("I100").Value = 17886.285
$ws.Range("J100").Value = 3949
$ws.Range("K100").Value = 17886.285
$ws.Range("L100").Value = 3949
$ws.Range("M100").Value = -17345.285
$ws.Range("N100").Value = -5031

# @@ -7211,22 +7211,22 @@  (sheet ALC)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2013.5416
$ws.Range("I132").Value = 2065.6956
$ws.Range("K132").Value = 6197.0868
$ws.Range("M132").Value = -3667.0868

# @@ -10647,22 +10647,22 @@  (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 76931190
$ws.Range("I61").Value = 76931190
$ws.Range("K61").Value = 76931190
$ws.Range("M61").Value = -76930978

# @@ -11931,25 +11931,25 @@  (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1868.0714
$ws.Range("I88").Value = 1375.6666
$ws.Range("J88").Value = 2237.375
$ws.Range("K88").Value = 1375.6666
$ws.Range("L88").Value = 2237.375
$ws.Range("M88").Value = -969.6666
$ws.Range("N88").Value = -3049.375

# @@ -12078,25 +12078,25 @@  (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1868.0714
$ws.Range("I91").Value = 1375.6666
$ws.Range("J91").Value = 2237.375
$ws.Range("K91").Value = 1375.6666
$ws.Range("L91").Value = 2237.375
$ws.Range("M91").Value = 28.33339999999998
$ws.Range("N91").Value = -5045.375

# @@ -14075,22 +14075,22 @@  (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2944950.8
$ws.Range("I132").Value = 3034085.5
$ws.Range("K132").Value = 9102256.5
$ws.Range("M132").Value = -9099726.5

# @@ -14173,22 +14173,22 @@  (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 208797.8
$ws.Range("J134").Value = 208797.8
$ws.Range("L134").Value = 208797.8
$ws.Range("N134").Value = -218937.8

# @@ -14271,22 +14271,22 @@  (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 76931190
$ws.Range("I136").Value = 76931190
$ws.Range("K136").Value = 230793570
$ws.Range("M136").Value = -230791020

# @@ -14412,22 +14412,22 @@  (sheet ARM)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 98485.336
$ws.Range("J139").Value = 98485.336
$ws.Range("L139").Value = 98485.336
$ws.Range("N139").Value = -108765.336

# @@ -14907,22 +14907,22 @@  (sheet BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 508.7143
$ws.Range("I7").Value = 528.2
$ws.Range("K7").Value = 528.2
$ws.Range("M7").Value = -415.2

# @@ -18727,25 +18727,22 @@  (sheet BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4118.8
$ws.Range("I86").Value = 4118.8
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4118.8
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2995.8
$ws.Range("N86").Value = ""

# @@ -18877,25 +18874,22 @@  (sheet BSM)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4118.8
$ws.Range("I89").Value = 4118.8
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 20594
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -14978
$ws.Range("N89").Value = ""

# @@ -21529,25 +21523,25 @@  (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 189.66667
$ws.Range("J2").Value = 352.5
$ws.Range("L2").Value = 352.5
$ws.Range("N2").Value = -578.5

# @@ -22968,19 +22962,25 @@  (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10648.743
$ws.Range("I31").Value = 7824.625
$ws.Range("J31").Value = 13026.947
$ws.Range("K31").Value = 7824.625
$ws.Range("L31").Value = 13026.947
$ws.Range("M31").Value = -7529.625
$ws.Range("N31").Value = -13616.947

# @@ -23115,19 +23115,25 @@  (sheet CRP)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10648.743
$ws.Range("I34").Value = 7824.625
$ws.Range("J34").Value = 13026.947
$ws.Range("K34").Value = 7824.625
$ws.Range("L34").Value = 13026.947
$ws.Range("M34").Value = -7622.625
$ws.Range("N34").Value = -13430.947

# @@ -28943,25 +28949,25 @@  (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 279.875
$ws.Range("J12").Value = 392.14285
$ws.Range("L12").Value = 1176.42855
$ws.Range("N12").Value = -1522.42855

# @@ -29047,22 +29053,22 @@  (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 366
$ws.Range("I14").Value = 366
$ws.Range("K14").Value = 1098
$ws.Range("M14").Value = -925

# @@ -34118,25 +34124,25 @@  (sheet CUL)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3393.375
$ws.Range("I117").Value = 496
$ws.Range("J117").Value = 6290.75
$ws.Range("K117").Value = 1488
$ws.Range("L117").Value = 18872.25
$ws.Range("M117").Value = 1954
$ws.Range("N117").Value = -25756.25

# @@ -35461,25 +35467,25 @@  (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 661.2222
$ws.Range("I2").Value = 1082.7
$ws.Range("J2").Value = 134.375
$ws.Range("K2").Value = 1082.7
$ws.Range("L2").Value = 134.375
$ws.Range("M2").Value = -969.7
$ws.Range("N2").Value = -360.375

# @@ -39605,22 +39611,19 @@  (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""

# @@ -39749,22 +39752,19 @@  (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""

# @@ -41146,22 +41146,19 @@  (sheet GSM)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""

# @@ -43332,22 +43329,22 @@  (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2681.3572
$ws.Range("I22").Value = 2584.2856
$ws.Range("K22").Value = 2584.2856
$ws.Range("M22").Value = -2289.2856

# @@ -43571,22 +43568,22 @@  (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2681.3572
$ws.Range("I27").Value = 2584.2856
$ws.Range("K27").Value = 2584.2856
$ws.Range("M27").Value = -2477.2856

# @@ -44493,22 +44490,22 @@  (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 947.25
$ws.Range("I46").Value = 947.25
$ws.Range("K46").Value = 947.25
$ws.Range("M46").Value = -759.25

# @@ -45568,25 +45565,25 @@  (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8335833
$ws.Range("I68").Value = 12501000
$ws.Range("J68").Value = 5499
$ws.Range("K68").Value = 12501000
$ws.Range("L68").Value = 5499
$ws.Range("M68").Value = -12500251
$ws.Range("N68").Value = -6997

# @@ -45718,25 +45715,25 @@  (sheet LTW)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 8335833
$ws.Range("I71").Value = 12501000
$ws.Range("J71").Value = 5499
$ws.Range("K71").Value = 62505000
$ws.Range("L71").Value = 27495
$ws.Range("M71").Value = -62501256
$ws.Range("N71").Value = -34983

# @@ -49243,20 +49240,23 @@  (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 29749
$ws.Range("I2").Value = 29000
$ws.Range("K2").Value = 29000
$ws.Range("M2").Value = -28888

# @@ -52120,22 +52120,19 @@  (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = ""

# @@ -52261,22 +52258,19 @@  (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = ""

# @@ -52448,19 +52442,22 @@  (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 271000
$ws.Range("J69").Value = 271000
$ws.Range("L69").Value = 271000
$ws.Range("N69").Value = -272498

# @@ -52586,19 +52583,22 @@  (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 271000
$ws.Range("J72").Value = 271000
$ws.Range("L72").Value = 813000
$ws.Range("N72").Value = -820488

# @@ -54556,22 +54556,19 @@  (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""

# @@ -55478,25 +55475,25 @@  (sheet WVR)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 41678750
$ws.Range("I132").Value = 125005890
$ws.Range("J132").Value = 15184.5
$ws.Range("K132").Value = 375017670
$ws.Range("L132").Value = 45553.5
$ws.Range("M132").Value = -375015140
$ws.Range("N132").Value = -50613.5
